$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, pushing existing rows 7:20 down to 8:21
$ws.Rows("7:7").Insert()

# Populate the new row 7 with the new record (same pattern as surrounding rows)
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C7").Value = "Los Lagos"
$ws.Range("D7").Value = 44715
$ws.Range("E7").Value = 10
$ws.Range("F7").Value = 100112012
$ws.Range("G7").Value = "Espinaca"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 30
$ws.Range("K7").Value = 11000
$ws.Range("L7").Value = 11000
$ws.Range("M7").Value = 11000
$ws.Range("N7").Value = '$/cuna 10 kilos'
$ws.Range("O7").Value = "Región Metropolitana"
$ws.Range("P7").Value = 1100
$ws.Range("Q7").Value = 10
$ws.Range("R7").Value = "Hortaliza"
